$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Random color code"
$ws.Range("C6").Value = "fucntion should creatse a random 4 color code"
$ws.Range("D6").Value = "fucntion does create a for color code"
$ws.Range("E6").Value = "PASS"

# Add new test case row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Phone test"
$ws.Range("C7").Value = "all above test should work on phone"
$ws.Range("D7").Value = "working as expected, althoug some lag"
$ws.Range("E7").Value = "PASS"

# Update the selected cell to reflect where the user left off
$ws.Range("E7").Select()
